$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new row 7
$ws.Rows("7:7").Insert()
$ws.Rows("7:7").RowHeight = 285

# Apply base fill (copy format from B2 which already uses the fillId4 orange fill)
$ws.Range("B2").Copy()
$ws.Range("A7:G7").PasteSpecial(-4122)

# Values
$ws.Range("A7").Value = 42767
$ws.Range("B7").Value = 0.99930555555555556
$ws.Range("C7").Value = "Divide and Conquer Soft Deadline"
$ws.Range("F7").Value = "FIRE_PLACEHOLDER"
$ws.Range("G7").Value = "TA_PLACEHOLDER"

# Fonts: body default (not bold) for A/B, bold for C,D,E,F,G (font1 = bold 11)
$ws.Range("A7:B7").Font.Bold = $false
$ws.Range("A7:B7").Font.Size = 11
$ws.Range("C7:G7").Font.Bold = $true
$ws.Range("C7:G7").Font.Size = 11

# Number formats
$ws.Range("A7").NumberFormat = "m/d/yyyy"
$ws.Range("B7").NumberFormat = "h:mm:ss AM/PM"

# Alignment
$ws.Range("A7").HorizontalAlignment = -4108
$ws.Range("D7").WrapText = $true
$ws.Range("F7:G7").WrapText = $true
$ws.Range("F7:G7").VerticalAlignment = -4160

Write-Host "done"
